$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: politeness_score (B12) was stored as text "3"; convert it to a
# genuine numeric value 3.
$ws.Range("B12").Value = 3

# Add a brand new row 13 with the additional annotation record.
$ws.Range("A13").Value = "Ying Tang"
# B13 keeps "3" as text (matches the original textual representation),
# so force a leading apostrophe to keep Excel from coercing it to a number.
$ws.Range("B13").Value = "'3"
$ws.Range("C13").Value = "Not too surprisingly"
$ws.Range("D13").Value = "CRT"
$ws.Range("E13").Value = "MET"
$ws.Range("F13").Value = "afe80f3f-3501-40b4-a3d0-1ad1f86c76ec"
$ws.Range("G13").Value = "r1BRfhiab_annotated.xlsx"
$ws.Range("H13").Value = "Not too surprisingly, the standard multiclass losses do not have the desired property, however approaches that reduce multi-class to binary classification at training time do, namely unnormalized models with penalized log Z (self-normalization), the NCE approach, as well as (the natural in the proposed setting) binary classification loss."
